$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.245.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6063"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07114"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2816"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07672"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.825"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6367"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001005"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.080.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.916"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.217.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.024"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.097"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1293"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.457"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.839"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.832"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.743"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6533"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.562"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.765"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.223.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01755"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.564"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.980.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.544"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.499"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05540"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.40%  "
